$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 495
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 3541.5715
$ws.Range("J19").Value = 4512.8887
$ws.Range("L19").Value = 4512.8887
$ws.Range("N19").Value = -4862.8887
$ws.Range("H28").Value = 1075.5714
$ws.Range("I28").Value = 495.83334
$ws.Range("J28").Value = 1510.375
$ws.Range("K28").Value = 495.83334
$ws.Range("L28").Value = 1510.375
$ws.Range("M28").Value = -10.83334000000002
$ws.Range("N28").Value = -2480.375
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2730
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -2064
$ws.Range("H76").Value = 20964.166
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 20964.166
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184
$ws.Range("H96").Value = 3069.6667
$ws.Range("I96").Value = 2667
$ws.Range("K96").Value = 8001
$ws.Range("M96").Value = -6628
$ws.Range("H99").Value = 873.75
$ws.Range("I99").Value = 798.5
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 2395.5
$ws.Range("L99").Value = 3750
$ws.Range("M99").Value = -897.5
$ws.Range("N99").Value = -6746
$ws.Range("H100").Value = 3763.7144
$ws.Range("I100").Value = 3561.75
$ws.Range("K100").Value = 3561.75
$ws.Range("M100").Value = -3020.75
$ws.Range("H107").Value = 295.52173
$ws.Range("I107").Value = 205.65
$ws.Range("K107").Value = 205.65
$ws.Range("M107").Value = 1714.35
$ws.Range("H125").Value = 2461.5
$ws.Range("I125").Value = 1644
$ws.Range("J125").Value = 2811.8572
$ws.Range("K125").Value = 14796
$ws.Range("L125").Value = 25306.7148
$ws.Range("M125").Value = -12336
$ws.Range("N125").Value = -30226.7148
$ws.Range("H132").Value = 166670240
$ws.Range("I132").Value = 166670240
$ws.Range("K132").Value = 500010720
$ws.Range("M132").Value = -500008190

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2112.7144
$ws.Range("I102").Value = 2181.5833
$ws.Range("K102").Value = 2181.5833
$ws.Range("M102").Value = -559.5832999999998
$ws.Range("H104").Value = 19369.666
$ws.Range("J104").Value = 19369.666
$ws.Range("L104").Value = 19369.666
$ws.Range("N104").Value = -26357.666
$ws.Range("H110").Value = 27428.428
$ws.Range("I110").Value = 45375.875
$ws.Range("K110").Value = 45375.875
$ws.Range("M110").Value = -43330.875
$ws.Range("H122").Value = 2395.6667
$ws.Range("I122").Value = 1871.5
$ws.Range("K122").Value = 5614.5
$ws.Range("M122").Value = -3164.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 381.5
$ws.Range("I22").Value = 381.5
$ws.Range("K22").Value = 381.5
$ws.Range("M22").Value = -208.5
$ws.Range("H134").Value = 2564.4634
$ws.Range("J134").Value = 2825.4443
$ws.Range("L134").Value = 8476.332900000001
$ws.Range("N134").Value = -13546.3329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5093.9614
$ws.Range("I31").Value = 3266.4119
$ws.Range("K31").Value = 3266.4119
$ws.Range("M31").Value = -2971.4119
$ws.Range("H34").Value = 5093.9614
$ws.Range("I34").Value = 3266.4119
$ws.Range("K34").Value = 3266.4119
$ws.Range("M34").Value = -3064.4119
$ws.Range("H62").Value = 4201
$ws.Range("J62").Value = 3668.3333
$ws.Range("L62").Value = 3668.3333
$ws.Range("N62").Value = -4916.3333
$ws.Range("H65").Value = 4201
$ws.Range("J65").Value = 3668.3333
$ws.Range("L65").Value = 18341.6665
$ws.Range("N65").Value = -24581.6665
$ws.Range("H86").Value = 6879.5
$ws.Range("I86").Value = 5507.4
$ws.Range("K86").Value = 5507.4
$ws.Range("M86").Value = -4384.4
$ws.Range("H89").Value = 6879.5
$ws.Range("I89").Value = 5507.4
$ws.Range("K89").Value = 27537
$ws.Range("M89").Value = -21921
$ws.Range("H99").Value = 17247.455
$ws.Range("I99").Value = 20187.938
$ws.Range("K99").Value = 20187.938
$ws.Range("M99").Value = -18689.938
$ws.Range("H126").Value = 17247.455
$ws.Range("I126").Value = 20187.938
$ws.Range("K126").Value = 60563.814
$ws.Range("M126").Value = -58093.814
$ws.Range("H131").Value = 39022.75
$ws.Range("J131").Value = 39022.75
$ws.Range("L131").Value = 39022.75
$ws.Range("N131").Value = -49102.75
$ws.Range("H132").Value = 5203
$ws.Range("I132").Value = 1906
$ws.Range("K132").Value = 5718
$ws.Range("M132").Value = -3188

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2338.4443
$ws.Range("J5").Value = 6380
$ws.Range("L5").Value = 19140
$ws.Range("N5").Value = -19364
$ws.Range("H107").Value = 1066.8096
$ws.Range("I107").Value = 606.625
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 1819.875
$ws.Range("L107").Value = 4050
$ws.Range("M107").Value = 100.125
$ws.Range("N107").Value = -7890
$ws.Range("H113").Value = 2187.0454
$ws.Range("J113").Value = 2911.7144
$ws.Range("L113").Value = 8735.143199999999
$ws.Range("N113").Value = -13075.1432
$ws.Range("H120").Value = 14999
$ws.Range("J120").Value = 14999
$ws.Range("L120").Value = 44997
$ws.Range("N120").Value = -54673
$ws.Range("H122").Value = 748.5
$ws.Range("I122").Value = 804
$ws.Range("K122").Value = 7236
$ws.Range("M122").Value = -4786
$ws.Range("H128").Value = 120931.664
$ws.Range("I128").Value = 120931.664
$ws.Range("K128").Value = 362794.992
$ws.Range("M128").Value = -357814.992
$ws.Range("H129").Value = 1500.871
$ws.Range("I129").Value = 558.55554
$ws.Range("J129").Value = 1886.3636
$ws.Range("K129").Value = 1675.66662
$ws.Range("L129").Value = 5659.0908
$ws.Range("M129").Value = 3324.33338
$ws.Range("N129").Value = -15659.0908
$ws.Range("H131").Value = 1364.2927
$ws.Range("J131").Value = 1557.75
$ws.Range("L131").Value = 4673.25
$ws.Range("N131").Value = -14753.25
$ws.Range("H133").Value = 19980
$ws.Range("J133").Value = 19980
$ws.Range("L133").Value = 59940
$ws.Range("N133").Value = -70060
$ws.Range("H135").Value = 2338.4443
$ws.Range("J135").Value = 6380
$ws.Range("L135").Value = 57420
$ws.Range("N135").Value = -62490

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3968.7354
$ws.Range("I113").Value = 3431.3684
$ws.Range("J113").Value = 4649.4
$ws.Range("K113").Value = 3431.3684
$ws.Range("L113").Value = 4649.4
$ws.Range("M113").Value = -1261.3684
$ws.Range("N113").Value = -8989.4
$ws.Range("H122").Value = 4533.3335
$ws.Range("I122").Value = 4400
$ws.Range("K122").Value = 13200
$ws.Range("M122").Value = -10750

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3816.4167
$ws.Range("I7").Value = 3149.6667
$ws.Range("J7").Value = 4483.1665
$ws.Range("K7").Value = 3149.6667
$ws.Range("L7").Value = 4483.1665
$ws.Range("M7").Value = -3037.6667
$ws.Range("N7").Value = -4707.1665
$ws.Range("H68").Value = 3374.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 3374.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 3891.878
$ws.Range("I122").Value = 2941.5386
$ws.Range("J122").Value = 5539.1333
$ws.Range("K122").Value = 8824.6158
$ws.Range("L122").Value = 16617.3999
$ws.Range("M122").Value = -6374.6158
$ws.Range("N122").Value = -21517.3999
$ws.Range("H126").Value = 3816.4167
$ws.Range("I126").Value = 3149.6667
$ws.Range("J126").Value = 4483.1665
$ws.Range("K126").Value = 9449.000100000001
$ws.Range("L126").Value = 13449.4995
$ws.Range("M126").Value = -6979.000100000001
$ws.Range("N126").Value = -18389.4995
$ws.Range("H131").Value = 84388
$ws.Range("J131").Value = 84388
$ws.Range("L131").Value = 84388
$ws.Range("N131").Value = -94468

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 67333730
$ws.Range("I100").Value = 72143140
$ws.Range("K100").Value = 144286280
$ws.Range("M100").Value = -144285739
$ws.Range("H126").Value = 11067.333
$ws.Range("I126").Value = 12851
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 38553
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -36083
$ws.Range("N126").Value = -27440
